# Update "想去人数" (want-to-go count) / "最低票价" (lowest price) figures
# across the workbook's sheets to match refreshed source data
# (gh-pages output regenerated at commit 456a3b4).

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsShow    = $wb.Worksheets.Item("演出")
$wsAll     = $wb.Worksheets.Item("全部类型")

# --- 展览 (Exhibitions) sheet ---
$wsExhibit.Range("F4").Value = 69
$wsExhibit.Range("F5").Value = 95
$wsExhibit.Range("F6").Value = 498
$wsExhibit.Range("F7").Value = 4856
$wsExhibit.Range("F8").Value = 4856
$wsExhibit.Range("F10").Value = 128
$wsExhibit.Range("F12").Value = 32
$wsExhibit.Range("F14").Value = 139
$wsExhibit.Range("F15").Value = 7843
$wsExhibit.Range("F16").Value = 259
$wsExhibit.Range("F17").Value = 133
$wsExhibit.Range("F19").Value = 557
$wsExhibit.Range("F20").Value = 1940
$wsExhibit.Range("F21").Value = 1941
$wsExhibit.Range("F22").Value = 6303
$wsExhibit.Range("F25").Value = 2100
$wsExhibit.Range("G26").Value = 55
$wsExhibit.Range("F27").Value = 2
$wsExhibit.Range("F28").Value = 6257
$wsExhibit.Range("F29").Value = 160
$wsExhibit.Range("F30").Value = 47
$wsExhibit.Range("F34").Value = 6627
$wsExhibit.Range("F36").Value = 31
$wsExhibit.Range("F38").Value = 3
$wsExhibit.Range("F39").Value = 4
$wsExhibit.Range("F42").Value = 26
$wsExhibit.Range("F44").Value = 2485
$wsExhibit.Range("F48").Value = 47
$wsExhibit.Range("F49").Value = 465
$wsExhibit.Range("F50").Value = 2177
$wsExhibit.Range("F51").Value = 56

# --- 演出 (Performances) sheet ---
$wsShow.Range("F7").Value = 42
$wsShow.Range("F9").Value = 45
$wsShow.Range("F10").Value = 9

# --- 全部类型 (All types) sheet ---
$wsAll.Range("F5").Value = 69
$wsAll.Range("F7").Value = 95
$wsAll.Range("F10").Value = 498
$wsAll.Range("F11").Value = 4856
$wsAll.Range("F12").Value = 4856
$wsAll.Range("F14").Value = 128
$wsAll.Range("F16").Value = 32
$wsAll.Range("F17").Value = 7843
$wsAll.Range("F18").Value = 259
$wsAll.Range("F19").Value = 133
$wsAll.Range("F20").Value = 557
$wsAll.Range("F21").Value = 1941
$wsAll.Range("F23").Value = 6303
$wsAll.Range("F25").Value = 42
$wsAll.Range("F27").Value = 2100
$wsAll.Range("G28").Value = 55
$wsAll.Range("F30").Value = 45
$wsAll.Range("F31").Value = 2
$wsAll.Range("F32").Value = 6257
$wsAll.Range("F33").Value = 160
$wsAll.Range("F34").Value = 9
$wsAll.Range("F35").Value = 47
$wsAll.Range("F37").Value = 6627
$wsAll.Range("F38").Value = 31
$wsAll.Range("F42").Value = 26
$wsAll.Range("F47").Value = 47
$wsAll.Range("F48").Value = 465
$wsAll.Range("F50").Value = 56
